$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.636.01'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.771.34'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '596.24'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '170.88'
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('D7').Value = '3.771.09'
$ws.Range('E7').Value = '  -1.76%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.166'
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.49'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.457'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000277'
$ws.Range('E13').Value = '  +6.86%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.82'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '4.400.29'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').Value = '3.772.21'
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '18.87'
$ws.Range('E17').Value = '  +3.49%  '
$ws.Range('D18').Value = '67.645.78'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.27'
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.60'
$ws.Range('E21').Value = '  -4.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '470.65'
$ws.Range('E22').Value = '  +0.94%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.723'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('E24').Value = '  -7.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '83.97'
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.20'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.44'
$ws.Range('E28').Value = '  +3.79%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.92'
$ws.Range('E30').Value = '  -1.77%  '
$ws.Range('D31').Value = '3.913.59'
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.74'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.26'
$ws.Range('E33').Value = '  -2.45%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '30.57'
$ws.Range('E34').Value = '  -2.18%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.18'
$ws.Range('E35').Value = '  -3.89%  '
$ws.Range('D36').Value = '3.731.88'
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.85'
$ws.Range('E37').Value = '  +6.59%  '
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('E41').Value = '  -1.54%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.316'
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.77'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '45.92'
$ws.Range('E47').Value = '  -2.13%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '399.07'
$ws.Range('E48').Value = '  -5.38%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.000271'
$ws.Range('E49').Value = '  -6.75%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '141.98'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0357'
$ws.Range('E51').Value = '  -0.21%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
